# Edit: retitle/repurpose the "Quantum Mechanics" article into an "Art" article,
# swap the byline/email, rewrite the body + summary paragraphs, and append a
# trailing empty paragraph, matching the target unified diff.

$d = $word.ActiveDocument

# Vertical-tab char == Word's internal <w:br/> line-break marker in Range.Text.
$vt = [char]0x0B

# ---------------------------------------------------------------------------
# 1) Title
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('Unveiling the Enigmatic Realm of Quantum Mechanics') | Out-Null
if ($rng.Find.Found) {
    $rng.Text = 'Unveiling the Realm of Arts: Exploring Art''s Transformative Power'
}

# ---------------------------------------------------------------------------
# 2) Byline: "Dr. Helen Reed" (3 runs) -> "Anya Patel"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('Dr. Helen Reed') | Out-Null
if ($rng.Find.Found) {
    $rng.Text = 'Anya Patel'
}

# ---------------------------------------------------------------------------
# 3) Email local-part + domain
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('helenreed48@iqresmail') | Out-Null
if ($rng.Find.Found) {
    $rng.Text = 'anyapatel2413@protonmail'
}

# ---------------------------------------------------------------------------
# 4) Body paragraph (the long multi-sentence paragraph with the two manual
#    double line-breaks). Replace the whole span so the new sentences that
#    the diff adds are included inline, in order.
# ---------------------------------------------------------------------------
$bodyOld = 'Quantum mechanics, the perplexing field of physics, has opened a portal into the enigmatic realm of subatomic particles, challenging our understanding of the universe. Its peculiar rules, often defying our intuition, offer tantalizing glimpses into the hidden workings of matter and energy.' + $vt + $vt + 'Delving into the quantum realm, physicists have uncovered phenomena such as wave-particle duality, where particles exhibit both wave-like and particle-like properties. They have observed particles behaving as if they were entangled, even across vast distances, defying classical notions of causality. The very concept of time and space appears to warp and twist, as particles seem to defy the constraints of locality and simultaneity.' + $vt + $vt + 'The implications of quantum mechanics extend far beyond the laboratory. With the potential to revolutionize fields from medicine to computing, quantum technology holds the promise of transformative applications. Quantum computers, harnessing the power of superposition and entanglement, could solve complex computations exponentially faster than classical computers. This would usher in a new era of scientific breakthroughs, drug discovery, materials engineering, and cryptographic security.'

$bodyNew = 'Art, a kaleidoscope of colors and expressions, has long served as a medium to interpret and understand our world. Its beauty and complexity hold sway over minds and souls, blurring the lines between creativity, appreciation, and perception. Through various civilizations and cultures, art has shaped the collective human experience in myriad ways, reflecting the times and emotions that mold us. Like a symphony of colors and strokes, art''s canvas captures the heartbeat of humanity''s journey.' + $vt + $vt + 'As we ponder the impact of art, we find that it possesses transformative qualities that transcend boundaries. With every brushstroke, melody penned, or sculpture carved, art has the power to evoke emotions, spark dialogue, and bridge cultural divides. It acts as a mirror, reflecting who we are, and as a window, allowing us to peer into the lives and experiences of others. Through art, we discover empathy and compassion as we connect with the human spirit in its shared experiences of joy, sorrow, and wonder.' + $vt + $vt + 'In the tapestry of human expression, art becomes a sanctuary of self-discovery and self-expression. It provides a platform for individuals to showcase their uniqueness, to communicate their thoughts and feelings in ways words often fail to capture. Whether it be the vibrant hues of a painting, the resonating chords of a song, or the poignant lines of a poem, art allows us to access our inner selves, to understand our deepest desires and aspirations. In this process of self-exploration, art provides healing, allowing us to process emotions and transcend life''s challenges.'

$rng = $d.Content
$rng.Find.Execute($bodyOld) | Out-Null
if ($rng.Find.Found) {
    $rng.Text = $bodyNew
} else {
    Write-Host 'WARNING: body paragraph text not found for replacement'
}

# ---------------------------------------------------------------------------
# 5) Summary paragraph
# ---------------------------------------------------------------------------
$summaryOld = 'Quantum mechanics has unveiled a hidden realm where particles exhibit enigmatic properties that confound our classical intuition. The field holds immense promise for groundbreaking technologies, including quantum computers capable of solving complex problems far beyond the reach of conventional machines. As we continue to explore the enigmatic quantum world, we embark on a journey to comprehend the fundamental nature of reality itself.'

$summaryNew = 'In the vast landscape of human endeavors, art stands as a beacon of creativity and transformation. Through its ability to transcend boundaries, evoke emotions, and foster self-expression, art serves as a window into the depths of humanity. Whether it be the majesty of a symphony or the simplicity of a child''s drawing, art has an unparalleled ability to capture the complexities of our existence. Art can inspire us, heal us, and bring us closer together. It is a testament to the power of imagination and the enduring beauty of the human spirit.'

$rng = $d.Content
$rng.Find.Execute($summaryOld) | Out-Null
if ($rng.Find.Found) {
    $rng.Text = $summaryNew
} else {
    Write-Host 'WARNING: summary paragraph text not found for replacement'
}

# ---------------------------------------------------------------------------
# 6) Trailing empty paragraph appended after the summary paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null

Write-Host 'Done.'
